# "working script, only from localhost and only for new user"
#
# The sheet gains a new leading "user_id" column (A), the old "username"/
# "password" columns slide over to B/C, and a new "company" column (D) is
# added. Row 2/3 ("script_user_1"/"script_user_2") each get a generated
# user1/user2 id in B plus their existing password moved to C, and a new
# "script" value in D.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Remember the current (pre-edit) text so we can relocate it ---------
$oldB1 = $ws.Range("B1").Value()   # "password" -> moves to C1
$oldB2 = $ws.Range("B2").Value()   # "password_1" -> moves to C2
$oldB3 = $ws.Range("B3").Value()   # "password_2" -> moves to C3

# --- Row 1 (header): give C1/D1 the same look as B1 ("username") --------
$ws.Range("B1").Copy($ws.Range("C1"))
$ws.Range("B1").Copy($ws.Range("D1"))

$ws.Range("A1").Value = "user_id"
$ws.Range("B1").Value = "username"
$ws.Range("C1").Value = $oldB1
$ws.Range("D1").Value = "company"

# --- Row 2: give C2/D2 the same look as B2 -------------------------------
$ws.Range("B2").Copy($ws.Range("C2"))
$ws.Range("B2").Copy($ws.Range("D2"))

$ws.Range("B2").Value = "user1"
$ws.Range("C2").Value = $oldB2
$ws.Range("D2").Value = "script"

# --- Row 3: give C3/D3 the same look as B3 -------------------------------
$ws.Range("B3").Copy($ws.Range("C3"))
$ws.Range("B3").Copy($ws.Range("D3"))

$ws.Range("B3").Value = "user2"
$ws.Range("C3").Value = $oldB3
$ws.Range("D3").Value = "script"

# --- Restore the active selection Excel had recorded after the edit -----
$ws.Range("E8").Select()
